# Update Round_1 (D), Round_2 (E), Round_3 (F) dice score values for rows 2-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = @(2, 9, 1)
    3  = @(3, 2, 4)
    4  = @(0, 10, 6)
    5  = @(4, 4, 5)
    6  = @(2, 5, 1)
    7  = @(9, 5, 7)
    8  = @(4, 7, 9)
    9  = @(1, 6, 10)
    10 = @(9, 2, 10)
    11 = @(10, 2, 4)
    12 = @(10, 7, 6)
    13 = @(2, 4, 4)
    14 = @(7, 6, 6)
    15 = @(1, 9, 5)
    16 = @(6, 2, 8)
    17 = @(2, 4, 8)
    18 = @(1, 1, 3)
    19 = @(5, 7, 4)
    20 = @(6, 4, 2)
    21 = @(2, 0, 3)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
}
